$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1551.2812
$ws.Range("J17").Value = 1551.2812
$ws.Range("L17").Value = 4653.8436
$ws.Range("N17").Value = -4989.8436
$ws.Range("H132").Value = 12680.844
$ws.Range("I132").Value = 3883
$ws.Range("J132").Value = 14711.115
$ws.Range("K132").Value = 11649
$ws.Range("L132").Value = 44133.345
$ws.Range("M132").Value = -9119
$ws.Range("N132").Value = -49193.345
$ws.Range("H138").Value = 6735.913
$ws.Range("J138").Value = 7678.8687
$ws.Range("L138").Value = 23036.6061
$ws.Range("N138").Value = -33316.6061

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3840.7258
$ws.Range("I32").Value = 2251.75
$ws.Range("J32").Value = 12103.4
$ws.Range("K32").Value = 2251.75
$ws.Range("L32").Value = 12103.4
$ws.Range("M32").Value = -1964.75
$ws.Range("N32").Value = -12677.4
$ws.Range("H45").Value = 5400
$ws.Range("J45").Value = 3100
$ws.Range("L45").Value = 3100
$ws.Range("N45").Value = -3854
$ws.Range("H132").Value = 36066.633
$ws.Range("I132").Value = 45482.46
$ws.Range("K132").Value = 136447.38
$ws.Range("M132").Value = -133917.38

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1553.3478
$ws.Range("I86").Value = 1489.5
$ws.Range("K86").Value = 1489.5
$ws.Range("M86").Value = -366.5
$ws.Range("H89").Value = 1553.3478
$ws.Range("I89").Value = 1489.5
$ws.Range("K89").Value = 7447.5
$ws.Range("M89").Value = -1831.5
$ws.Range("H107").Value = 1313.0454
$ws.Range("I107").Value = 1442
$ws.Range("J107").Value = 1223.7693
$ws.Range("K107").Value = 1442
$ws.Range("L107").Value = 1223.7693
$ws.Range("M107").Value = 478
$ws.Range("N107").Value = -5063.7693
$ws.Range("H134").Value = 3333.2856
$ws.Range("I134").Value = 2774.2727
$ws.Range("K134").Value = 8322.8181
$ws.Range("M134").Value = -5787.8181

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5580.7407
$ws.Range("I99").Value = 3759.2
$ws.Range("K99").Value = 3759.2
$ws.Range("M99").Value = -2261.2
$ws.Range("H105").Value = 1069.7307
$ws.Range("I105").Value = 969.8095
$ws.Range("J105").Value = 1489.4
$ws.Range("K105").Value = 969.8095
$ws.Range("L105").Value = 1489.4
$ws.Range("M105").Value = 777.1905
$ws.Range("N105").Value = -4983.4
$ws.Range("H122").Value = 2890.1853
$ws.Range("J122").Value = 5918.5
$ws.Range("L122").Value = 17755.5
$ws.Range("N122").Value = -22655.5
$ws.Range("H126").Value = 5580.7407
$ws.Range("I126").Value = 3759.2
$ws.Range("K126").Value = 11277.6
$ws.Range("M126").Value = -8807.599999999999
$ws.Range("H132").Value = 7938184.5
$ws.Range("I132").Value = 9525327
$ws.Range("K132").Value = 28575981
$ws.Range("M132").Value = -28573451
$ws.Range("H134").Value = 1810.7778
$ws.Range("I134").Value = 1857.4073
$ws.Range("K134").Value = 5572.2219
$ws.Range("M134").Value = -3037.2219

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 384
$ws.Range("J12").Value = 401.6
$ws.Range("L12").Value = 1204.8
$ws.Range("N12").Value = -1550.8
$ws.Range("H18").Value = 2845.3
$ws.Range("I18").Value = 1337.8
$ws.Range("J18").Value = 4352.8
$ws.Range("K18").Value = 4013.4
$ws.Range("L18").Value = 13058.4
$ws.Range("M18").Value = -3844.4
$ws.Range("N18").Value = -13396.4
$ws.Range("H69").Value = 4866.6665
$ws.Range("I69").Value = 600
$ws.Range("J69").Value = 7000
$ws.Range("K69").Value = 1800
$ws.Range("L69").Value = 21000
$ws.Range("M69").Value = -989
$ws.Range("N69").Value = -22622
$ws.Range("H72").Value = 4866.6665
$ws.Range("I72").Value = 600
$ws.Range("J72").Value = 7000
$ws.Range("K72").Value = 5400
$ws.Range("L72").Value = 63000
$ws.Range("M72").Value = -1344
$ws.Range("N72").Value = -71112
$ws.Range("H101").Value = 10998
$ws.Range("J101").Value = 10998
$ws.Range("L101").Value = 32994
$ws.Range("N101").Value = -37862
$ws.Range("H108").Value = 5820.5557
$ws.Range("I108").Value = 340.7143
$ws.Range("J108").Value = 25000
$ws.Range("K108").Value = 1022.1429
$ws.Range("L108").Value = 75000
$ws.Range("M108").Value = 1857.8571
$ws.Range("N108").Value = -80760
$ws.Range("H124").Value = 6800.3335
$ws.Range("I124").Value = 950.5
$ws.Range("J124").Value = 18500
$ws.Range("K124").Value = 2851.5
$ws.Range("L124").Value = 55500
$ws.Range("M124").Value = 2058.5
$ws.Range("N124").Value = -65320
$ws.Range("H138").Value = 3597.5334
$ws.Range("J138").Value = 3999
$ws.Range("L138").Value = 11997
$ws.Range("N138").Value = -22277
$ws.Range("H139").Value = 3176.842
$ws.Range("I139").Value = 1843.2307
$ws.Range("J139").Value = 6066.3335
$ws.Range("K139").Value = 5529.6921
$ws.Range("L139").Value = 18199.0005
$ws.Range("M139").Value = -389.6921000000002
$ws.Range("N139").Value = -28479.0005

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 23356.715
$ws.Range("J80").Value = 28544.908
$ws.Range("L80").Value = 28544.908
$ws.Range("N80").Value = -30540.908
$ws.Range("H83").Value = 23356.715
$ws.Range("J83").Value = 28544.908
$ws.Range("L83").Value = 142724.54
$ws.Range("N83").Value = -152708.54
$ws.Range("H126").Value = 4353
$ws.Range("I126").Value = 2924.8823
$ws.Range("J126").Value = 6220.5386
$ws.Range("K126").Value = 8774.6469
$ws.Range("L126").Value = 18661.6158
$ws.Range("M126").Value = -6304.6469
$ws.Range("N126").Value = -23601.6158

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6255.1816
$ws.Range("I96").Value = 5853.8
$ws.Range("K96").Value = 5853.8
$ws.Range("M96").Value = -4480.8
$ws.Range("H126").Value = 2908
$ws.Range("I126").Value = 2908
$ws.Range("K126").Value = 8724
$ws.Range("M126").Value = -6254
$ws.Range("H136").Value = 7103.09
$ws.Range("I136").Value = 2494.2727
$ws.Range("J136").Value = 9373.104499999999
$ws.Range("K136").Value = 7482.8181
$ws.Range("L136").Value = 28119.3135
$ws.Range("M136").Value = -4932.8181
$ws.Range("N136").Value = -33219.3135
